{"js": "// Prepend \"Design: \" to the start of each \"Answers\" paragraph in the\n// NB518 feedback table (one occurrence per table row that begins a\n// ListBullet answer paragraph).\nconst targets = [\n  \"Way too many systems but fortunately most systems same as in previous vessels.\",\n  \"More-less straight forward. Only Panama arrangement not fulfilled due to aft ship design.\",\n  \"I consider material handling quite smooth. Mostly because it was already third vessel.\",\n  \"Block manufacturing timetable and detail design areas could have been more in line to give extra time for design.\",\n  \"In my opinion communication was smooth between different parties.\",\n  \"We got us well employed by changing the hull structure, suppliers and sub-contractors in third vessel.\"\n];\n\nfor (const target of targets) {\n  const results = context.document.body.search(target, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find target text: \" + target);\n  }\n\n  // Insert the prefix right at the start of the matched range, without\n  // disturbing the rest of the run's text/formatting.\n  results.items[0].insertText(\"Design: \", Word.InsertLocation.start);\n}\n\nawait context.sync();\n", "ps1": "# Prepend \"Design: \" to the start of each \"Answers\" paragraph in the\n# NB518 feedback table (one occurrence per table row that begins a\n# ListBullet answer paragraph).\n$d = $word.ActiveDocument\n\n$targets = @(\n  \"Way too many systems but fortunately most systems same as in previous vessels.\",\n  \"More-less straight forward. Only Panama arrangement not fulfilled due to aft ship design.\",\n  \"I consider material handling quite smooth. Mostly because it was already third vessel.\",\n  \"Block manufacturing timetable and detail design areas could have been more in line to give extra time for design.\",\n  \"In my opinion communication was smooth between different parties.\",\n  \"We got us well employed by changing the hull structure, suppliers and sub-contractors in third vessel.\"\n)\n\nforeach ($t in $targets) {\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Text = $t\n  $find.Forward = $true\n  $find.Wrap = 0\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  if ($find.Execute()) {\n    $range.SetRange($range.Start, $range.Start)\n    $range.InsertBefore(\"Design: \")\n  }\n}\n"}
